# Re-sort the two-column dictionary (A2:B27) into Hungarian alphabetical
# order by the term in column A, keeping each term/definition pair
# together (commit: "Szótár abc sorrendben").
#
# The target order does not match a plain ordinal/ASCII sort of the
# accented Hungarian text, so rows are written out explicitly in the
# correct final order rather than relying on a generic Range.Sort.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = "Ajtó"
$ws.Range("B2").Value  = "Egyik szobából a másikba visz át. Lehet két- és egyirányú"

$ws.Range("A3").Value  = "Bénulás"
$ws.Range("B3").Value  = "Megegyezik a kábítással"

$ws.Range("A4").Value  = "Camembert"
$ws.Range("B4").Value  = "Olyan tárgy, ami használatkor gázt bocsát ki magából"

$ws.Range("A5").Value  = "Egyesülés"
$ws.Range("B5").Value  = "Az a folyamat mely során két szobából egy lesz"

$ws.Range("A6").Value  = "Gáz"
$ws.Range("B6").Value  = "Elkábítja a hallgatókat és oktatókat"

$ws.Range("A7").Value  = "Hallgató"
$ws.Range("B7").Value  = "Mérnökhallgató röviden."

$ws.Range("A8").Value  = "Kábítás"
$ws.Range("B8").Value  = "Az adott hallgató/oktató eldobja összes tárgyát, és nem vehet fel többet, amíg ez tart"

$ws.Range("A9").Value  = "Kibukás"
$ws.Range("B9").Value  = "A hallgató számára vége a játéknak"

$ws.Range("A10").Value = "Kör"
$ws.Range("B10").Value = "A játék idő mértékegysége, addig tart amíg minden karakter sorra nem kerül"

$ws.Range("A11").Value = "Labirintus"
$ws.Range("B11").Value = "Szobák sorozata"

$ws.Range("A12").Value = "Logarléc"
$ws.Range("B12").Value = "Az a tárgy amit ha felvesznek a hallgatók, megnyerhetik a játékot"

$ws.Range("A13").Value = "Maszk"
$ws.Range("B13").Value = "Az ezt viselő hallgatón nem hat a gáz"

$ws.Range("A14").Value = "Mérnökhallgató"
$ws.Range("B14").Value = "A játékos által írányított karakter"

$ws.Range("A15").Value = "Oktató"
$ws.Range("B15").Value = "A Mérnökhallgató ellensége, a játék által irányított karakter"

$ws.Range("A16").Value = "Osztódás"
$ws.Range("B16").Value = "Az a folyamat mely során egy szobából kettő lesz"

$ws.Range("A17").Value = "Összekapcsolás"
$ws.Range("B17").Value = "Tranzisztorok között végbe menő folyamat, ez után lehet az egyik tranzisztorral a másikhoz ugrani"

$ws.Range("A18").Value = "Rongy"
$ws.Range("B18").Value = "Olyan tárgy, ami adott ideig védettséget biztosít a vele egy szobában lévő oktatók ellen minden hallgatónak"

$ws.Range("A19").Value = "Söröspohár"
$ws.Range("B19").Value = "Olyan tárgy, ami védettséget biztosít a hallgatónak adott ideig"

$ws.Range("A20").Value = "Szoba"
$ws.Range("B20").Value = "Olyan hely, ahol mérnökhallgatók, oktatók és tárgyak lehetnek. A szobák között lehetnek ajtók"

$ws.Range("A21").Value = "Szomszédos"
$ws.Range("B21").Value = "Két szoba ilyen, ha van közöttük ajtó"

$ws.Range("A22").Value = "Tárgy"
$ws.Range("B22").Value = "Hallgatók és oktatók vehetik fel, rendelkezik egy képeséggel, több féle van"

$ws.Range("A23").Value = "Tárgy eldobása"
$ws.Range("B23").Value = "Az oktató/hallgató leteszi a szobába az adott tárgyat"

$ws.Range("A24").Value = "Tárgy felvétele"
$ws.Range("B24").Value = "Az oktató/hallgató magához veszi az adott tárgyat"

$ws.Range("A25").Value = "Tranzisztor"
$ws.Range("B25").Value = "Olyan tárgy, amely segítségével a hallgató egy szobából elugorhat egy másikba, még akkor is ha azok nem szomszédosak"

$ws.Range("A26").Value = "TVSZ"
$ws.Range("B26").Value = "Olyan tárgy, ami védettséget biztosít a hallgatónak 3 alkalommal"

$ws.Range("A27").Value = "Védettség"
$ws.Range("B27").Value = "Az oktató/gáz nem tudja bántani a hallgatót"
